# Story "SSDMS-36" ("date filter" story) was deleted because it was not
# a requirement (per commit message: Story ID SSDMS 36 Deleted).
#
# That story lived on row 39 of the "Product Backlog" sheet. Deleting the
# entire row shifts every following row up by one (rows 40-61 -> 39-60),
# and Excel automatically drops the now-orphaned shared strings
# ("SSDMS-36" / the date-filter story text) from the shared string table.
#
# Several of the remaining rows' "Owner" (column E) cells referenced
# people (Pooja Sharma, Alkesh, Deepak Kandpal, Ashish Mishra,
# Ekansh Kumar, Bhanu P Tiwari) whose names are also removed from the
# shared workbook along with the deleted story - those Owner cells are
# cleared out (left blank, keeping their existing style).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Backlog")

# Delete the whole row for SSDMS-36 - shifts rows 40:61 up to 39:60,
# fixes the sheet dimension, merged cell ranges, and compacts the shared
# string table for strings that were only referenced by this row.
$ws.Rows.Item(39).Delete()

# After the shift, clear the Owner (column E) value for every row whose
# owner was one of the six people removed along with this story.
$ownerRowsToClear = @(39, 40, 41, 42, 43, 44, 45, 55, 56, 57, 58, 59, 60)
foreach ($r in $ownerRowsToClear) {
    $ws.Cells.Item($r, 5).ClearContents()
}

# Leave the selection on the row that now occupies the old SSDMS-36 slot,
# matching the natural post-delete selection state.
$ws.Rows.Item(39).Select()
